$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.729.50"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.641.43"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0629"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0841"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.867.78"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "1.642.16"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("D17").Value = "26.738.86"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  +14.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("D34").Value = "1.292.02"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("E37").Value = "  -5.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.535"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.806"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "1.793.73"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.35%  "
